# Update gh-pages output data (regenerated scrape) across the four sheets:
#   1 = 展览 (Exhibitions), 2 = 演出 (Performances),
#   3 = 本地生活 (Local life), 4 = 全部类型 (All types, combined)
# Only column F ("想去人数" / interested-count) and, for two rows, column G
# ("最低票价" / min ticket price) change values to reflect refreshed data.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 8354
$ws.Range("F3").Value = 137
$ws.Range("F4").Value = 106
$ws.Range("F5").Value = 36112
$ws.Range("F7").Value = 613
$ws.Range("F8").Value = 737
$ws.Range("F9").Value = 468
$ws.Range("F10").Value = 151
$ws.Range("F13").Value = 71
$ws.Range("F14").Value = 645
$ws.Range("F15").Value = 465
$ws.Range("F16").Value = 27
$ws.Range("F17").Value = 591
$ws.Range("F18").Value = 165
$ws.Range("F19").Value = 436
$ws.Range("F20").Value = 431
$ws.Range("F21").Value = 1133
$ws.Range("F23").Value = 757
$ws.Range("F24").Value = 2412
$ws.Range("F25").Value = 906
$ws.Range("F26").Value = 518
$ws.Range("F27").Value = 83
$ws.Range("F28").Value = 1108
$ws.Range("F29").Value = 40
$ws.Range("F30").Value = 694
$ws.Range("G30").Value = 52.2
$ws.Range("F31").Value = 22
$ws.Range("F32").Value = 1110

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 303
$ws.Range("F4").Value = 359
$ws.Range("F5").Value = 323

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 577

# --- Sheet 4: 全部类型 (combined view of sheets 1-3) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 577
$ws.Range("F3").Value = 8354
$ws.Range("F4").Value = 137
$ws.Range("F5").Value = 106
$ws.Range("F6").Value = 303
$ws.Range("F7").Value = 36112
$ws.Range("F9").Value = 613
$ws.Range("F10").Value = 737
$ws.Range("F11").Value = 468
$ws.Range("F13").Value = 151
$ws.Range("F15").Value = 359
$ws.Range("F16").Value = 323
$ws.Range("F19").Value = 71
$ws.Range("F20").Value = 645
$ws.Range("F21").Value = 465
$ws.Range("F23").Value = 27
$ws.Range("F28").Value = 591
$ws.Range("F29").Value = 165
$ws.Range("F30").Value = 436
$ws.Range("F31").Value = 431
$ws.Range("F32").Value = 1133
$ws.Range("F34").Value = 757
$ws.Range("F35").Value = 2412
$ws.Range("F36").Value = 906
$ws.Range("F37").Value = 518
$ws.Range("F38").Value = 83
$ws.Range("F39").Value = 1108
$ws.Range("F40").Value = 40
$ws.Range("F42").Value = 694
$ws.Range("G42").Value = 52.2
$ws.Range("F43").Value = 22
$ws.Range("F44").Value = 1110
